{"js": "// The edit: the document's default (\"Primary\") header currently reads\n// \"Bild, S\u00e4ulen, Psalm23\". The author expanded this keyword list to\n// \"Psalm23, Fluss, Weg, Regenbogen, Bild, Hand, S\u00e4ulen\" and added a new\n// (empty) trailing paragraph to the header, styled the same as the\n// existing one (\"Kopfzeile\" / Header style).\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst section = sections.items[0];\nconst header = section.getHeader(\"Primary\");\nheader.load(\"paragraphs\");\nawait context.sync();\n\n// The header currently holds a single paragraph with the keyword list.\nconst firstParagraph = header.paragraphs.items[0];\nfirstParagraph.load(\"text\");\nawait context.sync();\n\n// Replace its text with the new, expanded keyword list (style/formatting\n// of the paragraph itself is left untouched).\nfirstParagraph.insertText(\n  \"Psalm23, Fluss, Weg, Regenbogen, Bild, Hand, S\u00e4ulen\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// Append a new, empty paragraph after it (same \"Kopfzeile\" header style).\nheader.insertParagraph(\"\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# The edit: the document's default (\"Primary\") header currently reads\n# \"Bild, S\u00e4ulen, Psalm23\". The author expanded this keyword list to\n# \"Psalm23, Fluss, Weg, Regenbogen, Bild, Hand, S\u00e4ulen\" and added a new\n# (empty) trailing paragraph to the header, styled the same as the\n# existing one (\"Kopfzeile\" / Header style).\n\n$d = $word.ActiveDocument\n$section = $d.Sections.Item(1)\n\n# wdHeaderFooterPrimary = 1 -> the default header (rId7 / header2.xml).\n$header = $section.Headers.Item(1)\n\n# Replace the keyword list text in the existing (only) paragraph.\n$header.Range.Text = \"Psalm23, Fluss, Weg, Regenbogen, Bild, Hand, S\u00e4ulen\"\n\n# Append a new, empty paragraph after it (inherits the \"Kopfzeile\" header style).\n$header.Range.InsertParagraphAfter()\n"}
